{"js": "// Helper: find the first (and expected only) match of `needle` in the\n// document body and replace it with `replacement`.\nasync function replaceOnce(context, needle, replacement, matchCase = true) {\n  const results = context.document.body.search(needle, { matchCase: matchCase });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + needle);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. Version number: 1.0 -> 1.2.5\nawait replaceOnce(context, \"1.0\", \"1.2.5\");\n\n// 2. Change type: Creation -> Update\nawait replaceOnce(context, \"Creation\", \"Update\");\n\n// 3. Date: 04/05/2023 -> 31/05/2023\nawait replaceOnce(context, \"04/05/2023\", \"31/05/2023\");\n\n// 4. Pre-condition text: fix accent + add trailing period\nawait replaceOnce(\n  context,\n  \"O usuario devidamente autenticado e na tela inicial do sistema\",\n  \"O usu\u00e1rio devidamente autenticado e na tela inicial do sistema.\"\n);\n\n// 5. Fix accents: numero -> n\u00famero, diarias -> di\u00e1rias (in step 2 text)\nawait replaceOnce(context, \"numero de diarias\", \"n\u00famero de di\u00e1rias\");\n\n// 6. Step 3: add period before the \"af[1,2,3,4]\" reference\nawait replaceOnce(\n  context,\n  \"3. Chefe Seleciona uma di\u00e1ria apta para pagamento af[1,2,3,4]\",\n  \"3. Chefe Seleciona uma di\u00e1ria apta para pagamento. af[1,2,3,4]\"\n);\n\n// 7. Step 4: add trailing period\nawait replaceOnce(\n  context,\n  \"4. System Destaca a di\u00e1ria selecionada \",\n  \"4. System Destaca a di\u00e1ria selecionada. \"\n);\n\n// 8. Step 2 (Detalhar Di\u00e1rias flow): add trailing period\nawait replaceOnce(\n  context,\n  \"2. System Apresenta a tela de Detalhar Di\u00e1rias \",\n  \"2. System Apresenta a tela de Detalhar Di\u00e1rias. \"\n);\n\n// 9. Step 2 (Atualiza lista): remove redundant \"o nome\" after \"onde\"\nawait replaceOnce(\n  context,\n  \"2. System Atualiza a lista de registros de solicita\u00e7\u00f5es, onde o nome dever\u00e1 constar o nome do usu\u00e1rio logado (que se atribuiu como respons\u00e1vel pela AP) no campo de atribui\u00e7\u00e3o (no caso de desatribui\u00e7\u00e3o, o nome dever\u00e1 ser removido). \",\n  \"2. System Atualiza a lista de registros de solicita\u00e7\u00f5es, onde dever\u00e1 constar o nome do usu\u00e1rio logado (que se atribuiu como respons\u00e1vel pela AP) no campo de atribui\u00e7\u00e3o (no caso de desatribui\u00e7\u00e3o, o nome dever\u00e1 ser removido). \"\n);\n\n// 10. Step 2 (Registrar Autoriza\u00e7\u00f5es de Pagamento): add trailing period\nawait replaceOnce(\n  context,\n  \"2. System Apresenta a tela de Registrar Autoriza\u00e7\u00f5es de Pagamento \",\n  \"2. System Apresenta a tela de Registrar Autoriza\u00e7\u00f5es de Pagamento. \"\n);\n", "ps1": "# Apply the \"1.2.4 -> 1.2.5\" version bump and accompanying minor text fixes.\n$d = $word.ActiveDocument\n\nfunction Replace-Once($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\n# 1. Version number: 1.0 -> 1.2.5\nReplace-Once \"1.0\" \"1.2.5\"\n\n# 2. Change type: Creation -> Update\nReplace-Once \"Creation\" \"Update\"\n\n# 3. Date: 04/05/2023 -> 31/05/2023\nReplace-Once \"04/05/2023\" \"31/05/2023\"\n\n# 4. Pre-condition text: fix accent + add trailing period\nReplace-Once \"O usuario devidamente autenticado e na tela inicial do sistema\" \"O usu\u00e1rio devidamente autenticado e na tela inicial do sistema.\"\n\n# 5. Fix accents: numero -> n\u00famero, diarias -> di\u00e1rias (in step 2 text)\nReplace-Once \"numero de diarias\" \"n\u00famero de di\u00e1rias\"\n\n# 6. Step 3: add period before the \"af[1,2,3,4]\" reference\nReplace-Once \"3. Chefe Seleciona uma di\u00e1ria apta para pagamento af[1,2,3,4]\" \"3. Chefe Seleciona uma di\u00e1ria apta para pagamento. af[1,2,3,4]\"\n\n# 7. Step 4: add trailing period\nReplace-Once \"4. System Destaca a di\u00e1ria selecionada \" \"4. System Destaca a di\u00e1ria selecionada. \"\n\n# 8. Step 2 (Detalhar Di\u00e1rias flow): add trailing period\nReplace-Once \"2. System Apresenta a tela de Detalhar Di\u00e1rias \" \"2. System Apresenta a tela de Detalhar Di\u00e1rias. \"\n\n# 9. Step 2 (Atualiza lista): remove redundant \"o nome\" after \"onde\"\nReplace-Once \"2. System Atualiza a lista de registros de solicita\u00e7\u00f5es, onde o nome dever\u00e1 constar o nome do usu\u00e1rio logado (que se atribuiu como respons\u00e1vel pela AP) no campo de atribui\u00e7\u00e3o (no caso de desatribui\u00e7\u00e3o, o nome dever\u00e1 ser removido). \" \"2. System Atualiza a lista de registros de solicita\u00e7\u00f5es, onde dever\u00e1 constar o nome do usu\u00e1rio logado (que se atribuiu como respons\u00e1vel pela AP) no campo de atribui\u00e7\u00e3o (no caso de desatribui\u00e7\u00e3o, o nome dever\u00e1 ser removido). \"\n\n# 10. Step 2 (Registrar Autoriza\u00e7\u00f5es de Pagamento): add trailing period\nReplace-Once \"2. System Apresenta a tela de Registrar Autoriza\u00e7\u00f5es de Pagamento \" \"2. System Apresenta a tela de Registrar Autoriza\u00e7\u00f5es de Pagamento. \"\n"}
